$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Коррекция бровей" (row 5) already has its name in A5.
# Fill in the rest of the procedures table for it: duration, cost and
# the Saturday (Сб, column I) availability window, matching the other rows.
$ws.Range("B5").NumberFormat = $ws.Range("B2").NumberFormat   # time format (h:mm)
$ws.Range("B5").Value = 40.0 / 1440.0                          # 0:40 duration
$ws.Range("C5").Value = 30                                     # cost
$ws.Range("I5").Value = "10:00-19:00"                           # Saturday slot

# Column I now holds data like its neighbours, so let it resize like they did.
[void]$ws.Columns("I").AutoFit()

# Leave the selection where the author ended up.
[void]$ws.Range("N13").Select()
